$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 256.4
$ws.Range("I33").Value = 196
$ws.Range("K33").Value = 196
$ws.Range("M33").Value = 33
$ws.Range("H39").Value = 234.26666
$ws.Range("I39").Value = 116.46154
$ws.Range("K39").Value = 349.38462
$ws.Range("M39").Value = -53.38461999999998
$ws.Range("H64").Value = 4327.909
$ws.Range("I64").Value = 3251
$ws.Range("J64").Value = 4943.2856
$ws.Range("K64").Value = 3251
$ws.Range("L64").Value = 4943.2856
$ws.Range("M64").Value = -3003
$ws.Range("N64").Value = -5439.2856
$ws.Range("H67").Value = 4327.909
$ws.Range("I67").Value = 3251
$ws.Range("J67").Value = 4943.2856
$ws.Range("K67").Value = 3251
$ws.Range("L67").Value = 4943.2856
$ws.Range("M67").Value = -2393
$ws.Range("N67").Value = -6659.2856
$ws.Range("H80").Value = 1233.2222
$ws.Range("I80").Value = 1399.8572
$ws.Range("J80").Value = 650
$ws.Range("K80").Value = 4199.571599999999
$ws.Range("L80").Value = 1950
$ws.Range("M80").Value = -3201.571599999999
$ws.Range("N80").Value = -3946
$ws.Range("H83").Value = 1233.2222
$ws.Range("I83").Value = 1399.8572
$ws.Range("J83").Value = 650
$ws.Range("K83").Value = 12598.7148
$ws.Range("L83").Value = 5850
$ws.Range("M83").Value = -7606.7148
$ws.Range("N83").Value = -15834
$ws.Range("H103").Value = 351.72726
$ws.Range("I103").Value = 332.16666
$ws.Range("J103").Value = 375.2
$ws.Range("K103").Value = 996.4999799999999
$ws.Range("L103").Value = 1125.6
$ws.Range("M103").Value = -410.4999799999999
$ws.Range("N103").Value = -2297.6
$ws.Range("H107").Value = 917.4286
$ws.Range("I107").Value = 917.4286
$ws.Range("K107").Value = 917.4286
$ws.Range("M107").Value = 1002.5714
$ws.Range("H111").Value = 2822.4614
$ws.Range("I111").Value = 1189.1111
$ws.Range("J111").Value = 6497.5
$ws.Range("K111").Value = 3567.3333
$ws.Range("L111").Value = 19492.5
$ws.Range("M111").Value = -500.3333000000002
$ws.Range("N111").Value = -25626.5
$ws.Range("H112").Value = 1577.9474
$ws.Range("J112").Value = 1780.125
$ws.Range("L112").Value = 5340.375
$ws.Range("N112").Value = -7556.375
$ws.Range("H113").Value = 2713.5715
$ws.Range("I113").Value = 2799
$ws.Range("K113").Value = 2799
$ws.Range("M113").Value = 455
$ws.Range("H116").Value = 3873
$ws.Range("I116").Value = 3746.5
$ws.Range("K116").Value = 3746.5
$ws.Range("M116").Value = -304.5
$ws.Range("H132").Value = 916.76086
$ws.Range("I132").Value = 871.44183
$ws.Range("K132").Value = 2614.32549
$ws.Range("M132").Value = -84.32549000000017
$ws.Range("H137").Value = 1957.2
$ws.Range("I137").Value = 2143.25
$ws.Range("J137").Value = 1833.1666
$ws.Range("K137").Value = 6429.75
$ws.Range("L137").Value = 5499.4998
$ws.Range("M137").Value = -3879.75
$ws.Range("N137").Value = -10599.4998
$ws.Range("H138").Value = 5342.091
$ws.Range("I138").Value = 2819.9023
$ws.Range("K138").Value = 8459.706900000001
$ws.Range("M138").Value = -3319.706900000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5360.837
$ws.Range("I32").Value = 2802.1428
$ws.Range("K32").Value = 2802.1428
$ws.Range("M32").Value = -2515.1428
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()
$ws.Range("H74").Value = 1855
$ws.Range("I74").Value = 1200
$ws.Range("J74").Value = 2291.6667
$ws.Range("K74").Value = 1200
$ws.Range("L74").Value = 2291.6667
$ws.Range("M74").Value = -326
$ws.Range("N74").Value = -4039.6667
$ws.Range("H76").Value = 13333
$ws.Range("J76").Value = 13333
$ws.Range("L76").Value = 13333
$ws.Range("N76").Value = -14009
$ws.Range("H77").Value = 1855
$ws.Range("I77").Value = 1200
$ws.Range("J77").Value = 2291.6667
$ws.Range("K77").Value = 6000
$ws.Range("L77").Value = 11458.3335
$ws.Range("M77").Value = -1632
$ws.Range("N77").Value = -20194.3335
$ws.Range("H79").Value = 13333
$ws.Range("J79").Value = 13333
$ws.Range("L79").Value = 13333
$ws.Range("N79").Value = -15673

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 267.5
$ws.Range("I5").Value = 316.66666
$ws.Range("J5").Value = 120
$ws.Range("K5").Value = 316.66666
$ws.Range("L5").Value = 120
$ws.Range("M5").Value = -203.66666
$ws.Range("N5").Value = -346
$ws.Range("H20").Value = 5887.857
$ws.Range("I20").Value = 4751.5
$ws.Range("K20").Value = 4751.5
$ws.Range("M20").Value = -4504.5
$ws.Range("H105").Value = 4444.778
$ws.Range("I105").Value = 4444.778
$ws.Range("K105").Value = 4444.778
$ws.Range("M105").Value = -2697.778

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3768.4482
$ws.Range("I31").Value = 3125.375
$ws.Range("J31").Value = 6855.2
$ws.Range("K31").Value = 3125.375
$ws.Range("L31").Value = 6855.2
$ws.Range("M31").Value = -2830.375
$ws.Range("N31").Value = -7445.2
$ws.Range("H34").Value = 3768.4482
$ws.Range("I34").Value = 3125.375
$ws.Range("J34").Value = 6855.2
$ws.Range("K34").Value = 3125.375
$ws.Range("L34").Value = 6855.2
$ws.Range("M34").Value = -2923.375
$ws.Range("N34").Value = -7259.2
$ws.Range("H60").Value = 16069.833
$ws.Range("I60").Value = 12985.363
$ws.Range("J60").Value = 49999
$ws.Range("K60").Value = 12985.363
$ws.Range("L60").Value = 49999
$ws.Range("M60").Value = -12474.363
$ws.Range("N60").Value = -51021
$ws.Range("H107").Value = 2177.037
$ws.Range("I107").Value = 2522.7144
$ws.Range("K107").Value = 2522.7144
$ws.Range("M107").Value = -602.7143999999998
$ws.Range("H134").Value = 2427.889
$ws.Range("I134").Value = 1896.4166
$ws.Range("J134").Value = 3490.8333
$ws.Range("K134").Value = 5689.2498
$ws.Range("L134").Value = 10472.4999
$ws.Range("M134").Value = -3154.2498
$ws.Range("N134").Value = -15542.4999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H28").Value = 1465
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").ClearContents()
$ws.Range("H107").Value = 1968.3334
$ws.Range("J107").Value = 1968.3334
$ws.Range("L107").Value = 5905.0002
$ws.Range("N107").Value = -9745.0002
$ws.Range("H118").Value = 585.3
$ws.Range("I118").Value = 585.3
$ws.Range("K118").Value = 1755.9
$ws.Range("M118").Value = -512.8999999999999
$ws.Range("H125").Value = 1000
$ws.Range("I125").Value = 1000
$ws.Range("K125").Value = 3000
$ws.Range("M125").Value = 1920

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 8613.25
$ws.Range("I80").Value = 5380.4
$ws.Range("K80").Value = 5380.4
$ws.Range("M80").Value = -4382.4
$ws.Range("H83").Value = 8613.25
$ws.Range("I83").Value = 5380.4
$ws.Range("K83").Value = 26902
$ws.Range("M83").Value = -21910
$ws.Range("H113").Value = 1999
$ws.Range("I113").Value = 1999
$ws.Range("K113").Value = 1999
$ws.Range("M113").Value = 171

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 6210.162
$ws.Range("I22").Value = 5018.4443
$ws.Range("J22").Value = 7339.1577
$ws.Range("K22").Value = 5018.4443
$ws.Range("L22").Value = 7339.1577
$ws.Range("M22").Value = -4723.4443
$ws.Range("N22").Value = -7929.1577
$ws.Range("H27").Value = 6210.162
$ws.Range("I27").Value = 5018.4443
$ws.Range("J27").Value = 7339.1577
$ws.Range("K27").Value = 5018.4443
$ws.Range("L27").Value = 7339.1577
$ws.Range("M27").Value = -4911.4443
$ws.Range("N27").Value = -7553.1577
$ws.Range("H40").Value = 2853.3076
$ws.Range("I40").Value = 2663
$ws.Range("K40").Value = 2663
$ws.Range("M40").Value = -2527
$ws.Range("H61").Value = 1313.8334
$ws.Range("I61").Value = 1313.8334
$ws.Range("K61").Value = 1313.8334
$ws.Range("M61").Value = -1111.8334
$ws.Range("H113").Value = 1313.8334
$ws.Range("I113").Value = 1313.8334
$ws.Range("K113").Value = 1313.8334
$ws.Range("M113").Value = 856.1666
$ws.Range("H136").Value = 4223.625
$ws.Range("I136").Value = 4112.7144
$ws.Range("K136").Value = 12338.1432
$ws.Range("M136").Value = -9788.143199999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 3000000
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()
$ws.Range("H81").Value = 9264.190000000001
$ws.Range("I81").Value = 2649.7144
$ws.Range("K81").Value = 5299.4288
$ws.Range("M81").Value = -4238.4288
$ws.Range("H84").Value = 9264.190000000001
$ws.Range("I84").Value = 2649.7144
$ws.Range("K84").Value = 26497.144
$ws.Range("M84").Value = -21193.144
$ws.Range("H126").Value = 3488.75
$ws.Range("I126").Value = 1935
$ws.Range("K126").Value = 5805
$ws.Range("M126").Value = -3335
$ws.Range("H132").Value = 44053.684
$ws.Range("I132").Value = 59636.688
$ws.Range("K132").Value = 178910.064
$ws.Range("M132").Value = -176380.064
